$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 223, shifting existing rows 223:250 down to 224:251
$ws.Rows("223:223").Insert()

# Populate the newly inserted row 223 with the new record
$r = 223
$ws.Cells.Item($r, 1).Value = 11
$ws.Cells.Item($r, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item($r, 3).Value = "Bíobío"
$ws.Cells.Item($r, 4).Value = 44776
$ws.Cells.Item($r, 4).NumberFormat = $ws.Cells.Item(224, 4).NumberFormat
$ws.Cells.Item($r, 5).Value = 8
$ws.Cells.Item($r, 6).Value = 100114013
$ws.Cells.Item($r, 7).Value = "Zanahoria"
$ws.Cells.Item($r, 8).Value = "Sin especificar"
$ws.Cells.Item($r, 9).Value = "Primera"
$ws.Cells.Item($r, 10).Value = 350
$ws.Cells.Item($r, 11).Value = 8000
$ws.Cells.Item($r, 12).Value = 9000
$ws.Cells.Item($r, 13).Value = 8429
$ws.Cells.Item($r, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item($r, 15).Value = "Chillán"
$ws.Cells.Item($r, 16).Value = 421
$ws.Cells.Item($r, 17).Value = 20
$ws.Cells.Item($r, 18).Value = "Hortaliza"
